# Auto-generated edit script for herbariumList.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells E1 / F1 ---
$ws.Cells.Item(1, 5).Value = "Dry Weight"
$ws.Cells.Item(1, 6).Value = "Solvent volume for 10000 ppm (µL)"

# Bold the "µL)" portion of F1 (rich-text run), matching the source formatting
$f1 = $ws.Range("F1")
$f1.Characters(31, 3).Font.Bold = $true

# --- Column widths for new columns E (15.140625) and F (13.28515625) ---
$ws.Columns.Item(5).ColumnWidth = 14.25
$ws.Columns.Item(6).ColumnWidth = 12.42

# --- Clear the "text" style (s=1, right-aligned @ format) from column A so new numbers store as real numbers ---
$ws.Range("A46:A93").Style = "Normal"

# --- Apply vertical-center style (matches existing B/C data cells, style index 3) to new B/C cells ---
$ws.Range("B46:C103").VerticalAlignment = -4108

# row 46
$ws.Cells.Item(46, 1).Value = 31
$ws.Cells.Item(46, 2).Value = "NYBG03"
$ws.Cells.Item(46, 3).Value = "S. angustifolia ssp. micranthum"
$ws.Cells.Item(46, 5).Value = 0.0073
$ws.Cells.Item(46, 6).Formula = "=E46*1000/10*1000"

# row 47
$ws.Cells.Item(47, 1).Value = 32
$ws.Cells.Item(47, 2).Value = "NYBG04"
$ws.Cells.Item(47, 3).Value = "S. angustifolia ssp. angustifolia"
$ws.Cells.Item(47, 5).Value = 0.0153
$ws.Cells.Item(47, 6).Formula = "=E47*1000/10*1000"

# row 48
$ws.Cells.Item(48, 1).Value = 33
$ws.Cells.Item(48, 2).Value = "NYBG05"
$ws.Cells.Item(48, 3).Value = "S. antirrhinoides"
$ws.Cells.Item(48, 5).Value = 0.0091
$ws.Cells.Item(48, 6).Formula = "=E48*1000/10*1000"

# row 49
$ws.Cells.Item(49, 1).Value = 34
$ws.Cells.Item(49, 2).Value = "NYBG07"
$ws.Cells.Item(49, 3).Value = "S. brittonii"
$ws.Cells.Item(49, 5).Value = 0.0135
$ws.Cells.Item(49, 6).Formula = "=E49*1000/10*1000"

# row 50
$ws.Cells.Item(50, 1).Value = 35
$ws.Cells.Item(50, 2).Value = "NYBG08"
$ws.Cells.Item(50, 3).Value = "S. californica"
$ws.Cells.Item(50, 5).Value = 0.0047
$ws.Cells.Item(50, 6).Formula = "=E50*1000/10*1000"

# row 51
$ws.Cells.Item(51, 1).Value = 36
$ws.Cells.Item(51, 2).Value = "NYBG09"
$ws.Cells.Item(51, 3).Value = "S. coccinea"
$ws.Cells.Item(51, 5).Value = 0.0051
$ws.Cells.Item(51, 6).Formula = "=E51*1000/10*1000"

# row 52
$ws.Cells.Item(52, 1).Value = 37
$ws.Cells.Item(52, 2).Value = "NYBG10"
$ws.Cells.Item(52, 3).Value = "S. drumondii"
$ws.Cells.Item(52, 5).Value = 0.021
$ws.Cells.Item(52, 6).Formula = "=E52*1000/10*1000"

# row 53
$ws.Cells.Item(53, 1).Value = 38
$ws.Cells.Item(53, 2).Value = "NYBG13"
$ws.Cells.Item(53, 3).Value = "S. glabriuscula"
$ws.Cells.Item(53, 5).Value = 0.0088
$ws.Cells.Item(53, 6).Formula = "=E53*1000/10*1000"

# row 54
$ws.Cells.Item(54, 1).Value = 39
$ws.Cells.Item(54, 2).Value = "NYBG14"
$ws.Cells.Item(54, 3).Value = "S. havanensis"
$ws.Cells.Item(54, 5).Value = 0.0159
$ws.Cells.Item(54, 6).Formula = "=E54*1000/10*1000"

# row 55
$ws.Cells.Item(55, 1).Value = 40
$ws.Cells.Item(55, 2).Value = "NYBG15"
$ws.Cells.Item(55, 3).Value = "S. holmgrenierum"
$ws.Cells.Item(55, 5).Value = 0.0129
$ws.Cells.Item(55, 6).Formula = "=E55*1000/10*1000"

# row 56
$ws.Cells.Item(56, 1).Value = 41
$ws.Cells.Item(56, 2).Value = "NYBG17"
$ws.Cells.Item(56, 3).Value = "S. incana"
$ws.Cells.Item(56, 5).Value = 0.006
$ws.Cells.Item(56, 6).Formula = "=E56*1000/10*1000"

# row 57
$ws.Cells.Item(57, 1).Value = 42
$ws.Cells.Item(57, 2).Value = "NYBG19"
$ws.Cells.Item(57, 3).Value = "S. leonardii"
$ws.Cells.Item(57, 5).Value = 0.009
$ws.Cells.Item(57, 6).Formula = "=E57*1000/10*1000"

# row 58
$ws.Cells.Item(58, 1).Value = 43
$ws.Cells.Item(58, 2).Value = "NYBG22"
$ws.Cells.Item(58, 3).Value = "S. multiglandulosa"
$ws.Cells.Item(58, 5).Value = 0.0068
$ws.Cells.Item(58, 6).Formula = "=E58*1000/10*1000"

# row 59
$ws.Cells.Item(59, 1).Value = 44
$ws.Cells.Item(59, 2).Value = "NYBG23"
$ws.Cells.Item(59, 3).Value = "S. muriculata"
$ws.Cells.Item(59, 5).Value = 0.0138
$ws.Cells.Item(59, 6).Formula = "=E59*1000/10*1000"

# row 60
$ws.Cells.Item(60, 1).Value = 45
$ws.Cells.Item(60, 2).Value = "NYBG25"
$ws.Cells.Item(60, 3).Value = "S. sapphirina"
$ws.Cells.Item(60, 5).Value = 0.0119
$ws.Cells.Item(60, 6).Formula = "=E60*1000/10*1000"

# row 61
$ws.Cells.Item(61, 1).Value = 46
$ws.Cells.Item(61, 2).Value = "NYBG29"
$ws.Cells.Item(61, 3).Value = "S. pseudoserrata"
$ws.Cells.Item(61, 5).Value = 0.0057
$ws.Cells.Item(61, 6).Formula = "=E61*1000/10*1000"

# row 62
$ws.Cells.Item(62, 1).Value = 47
$ws.Cells.Item(62, 2).Value = "NYBG30"
$ws.Cells.Item(62, 3).Value = "S. racemosa"
$ws.Cells.Item(62, 5).Value = 0.0103
$ws.Cells.Item(62, 6).Formula = "=E62*1000/10*1000"

# row 63
$ws.Cells.Item(63, 1).Value = 48
$ws.Cells.Item(63, 2).Value = "NYBG32"
$ws.Cells.Item(63, 3).Value = "S. sapphirina"
$ws.Cells.Item(63, 5).Value = 0.0205
$ws.Cells.Item(63, 6).Formula = "=E63*1000/10*1000"

# row 64
$ws.Cells.Item(64, 1).Value = 49
$ws.Cells.Item(64, 2).Value = "NYBG38"
$ws.Cells.Item(64, 3).Value = "S. blepharophylla"
$ws.Cells.Item(64, 5).Value = 0.0096
$ws.Cells.Item(64, 6).Formula = "=E64*1000/10*1000"

# row 65
$ws.Cells.Item(65, 1).Value = 50
$ws.Cells.Item(65, 2).Value = "NYBG42"
$ws.Cells.Item(65, 3).Value = "S. hispidula"
$ws.Cells.Item(65, 5).Value = 0.0062
$ws.Cells.Item(65, 6).Formula = "=E65*1000/10*1000"

# row 66
$ws.Cells.Item(66, 1).Value = 51
$ws.Cells.Item(66, 2).Value = "NYBG44"
$ws.Cells.Item(66, 3).Value = "S. guatemalensis"
$ws.Cells.Item(66, 5).Value = 0.0112
$ws.Cells.Item(66, 6).Formula = "=E66*1000/10*1000"

# row 67
$ws.Cells.Item(67, 1).Value = 52
$ws.Cells.Item(67, 2).Value = "NYBG48"
$ws.Cells.Item(67, 3).Value = "S. lutea"
$ws.Cells.Item(67, 5).Value = 0.0225
$ws.Cells.Item(67, 6).Formula = "=E67*1000/10*1000"

# row 68
$ws.Cells.Item(68, 1).Value = 53
$ws.Cells.Item(68, 2).Value = "NYBG50"
$ws.Cells.Item(68, 3).Value = "S. purpurascens"
$ws.Cells.Item(68, 5).Value = 0.0138
$ws.Cells.Item(68, 6).Formula = "=E68*1000/10*1000"

# row 69
$ws.Cells.Item(69, 1).Value = 54
$ws.Cells.Item(69, 2).Value = "NYBG51"
$ws.Cells.Item(69, 3).Value = "S. seleriana"
$ws.Cells.Item(69, 5).Value = 0.0164
$ws.Cells.Item(69, 6).Formula = "=E69*1000/10*1000"

# row 70
$ws.Cells.Item(70, 1).Value = 55
$ws.Cells.Item(70, 2).Value = "NYBG55"
$ws.Cells.Item(70, 3).Value = "S. suffrutscens"
$ws.Cells.Item(70, 5).Value = 0.0199
$ws.Cells.Item(70, 6).Formula = "=E70*1000/10*1000"

# row 71
$ws.Cells.Item(71, 1).Value = 56
$ws.Cells.Item(71, 2).Value = "NYBG56"
$ws.Cells.Item(71, 3).Value = "S. przewalskii"
$ws.Cells.Item(71, 5).Value = 0.014
$ws.Cells.Item(71, 6).Formula = "=E71*1000/10*1000"

# row 72
$ws.Cells.Item(72, 1).Value = 57
$ws.Cells.Item(72, 2).Value = "NYBG57"
$ws.Cells.Item(72, 3).Value = "S. scordiifolia"
$ws.Cells.Item(72, 5).Value = 0.0287
$ws.Cells.Item(72, 6).Formula = "=E72*1000/10*1000"

# row 73
$ws.Cells.Item(73, 1).Value = 58
$ws.Cells.Item(73, 2).Value = "NYBG58"
$ws.Cells.Item(73, 3).Value = "S. discolor"
$ws.Cells.Item(73, 5).Value = 0.0064
$ws.Cells.Item(73, 6).Formula = "=E73*1000/10*1000"

# row 74
$ws.Cells.Item(74, 1).Value = 59
$ws.Cells.Item(74, 2).Value = "NYBG60"
$ws.Cells.Item(74, 3).Value = "S. multicularis"
$ws.Cells.Item(74, 5).Value = 0.0825
$ws.Cells.Item(74, 6).Formula = "=E74*1000/10*1000"

# row 75
$ws.Cells.Item(75, 1).Value = 60
$ws.Cells.Item(75, 2).Value = "NYBG61"
$ws.Cells.Item(75, 3).Value = "S. oblonga"
$ws.Cells.Item(75, 5).Value = 0.0123
$ws.Cells.Item(75, 6).Formula = "=E75*1000/10*1000"

# row 76
$ws.Cells.Item(76, 1).Value = 61
$ws.Cells.Item(76, 2).Value = "NYBG62"
$ws.Cells.Item(76, 3).Value = "S. heterophylla"

# row 77
$ws.Cells.Item(77, 1).Value = 62
$ws.Cells.Item(77, 2).Value = "NYBG63"
$ws.Cells.Item(77, 3).Value = "S. heydei"

# row 78
$ws.Cells.Item(78, 1).Value = 63
$ws.Cells.Item(78, 2).Value = "NYBG64"
$ws.Cells.Item(78, 3).Value = "S. javanica"

# row 79
$ws.Cells.Item(79, 1).Value = 64
$ws.Cells.Item(79, 2).Value = "NYBG65"
$ws.Cells.Item(79, 3).Value = "S. pinnatifida"

# row 80
$ws.Cells.Item(80, 1).Value = 65
$ws.Cells.Item(80, 2).Value = "NYBG66"
$ws.Cells.Item(80, 3).Value = "S. prostrata"

# row 81
$ws.Cells.Item(81, 1).Value = 66
$ws.Cells.Item(81, 2).Value = "NYBG67"
$ws.Cells.Item(81, 3).Value = "S. peregrina"

# row 82
$ws.Cells.Item(82, 1).Value = 67
$ws.Cells.Item(82, 2).Value = "NYBG68"
$ws.Cells.Item(82, 3).Value = "S. relenorskyi"

# row 83
$ws.Cells.Item(83, 1).Value = 68
$ws.Cells.Item(83, 2).Value = "NYBG69"
$ws.Cells.Item(83, 3).Value = "S. angulosa"

# row 84
$ws.Cells.Item(84, 1).Value = 69
$ws.Cells.Item(84, 2).Value = "NYBG71"
$ws.Cells.Item(84, 3).Value = "S. paucifolia"

# row 85
$ws.Cells.Item(85, 1).Value = 70
$ws.Cells.Item(85, 2).Value = "NYBG72"
$ws.Cells.Item(85, 3).Value = "S. sumatrana"

# row 86
$ws.Cells.Item(86, 1).Value = 71
$ws.Cells.Item(86, 2).Value = "NYBG73"
$ws.Cells.Item(86, 3).Value = "S. javalambrensis"

# row 87
$ws.Cells.Item(87, 1).Value = 72
$ws.Cells.Item(87, 2).Value = "NYBG74"
$ws.Cells.Item(87, 3).Value = "S. leptosiplonsipkon"

# row 88
$ws.Cells.Item(88, 1).Value = 73
$ws.Cells.Item(88, 2).Value = "FLAS01"
$ws.Cells.Item(88, 3).Value = "S. arenicola"

# row 89
$ws.Cells.Item(89, 1).Value = 74
$ws.Cells.Item(89, 2).Value = "FLAS02"
$ws.Cells.Item(89, 3).Value = "S. arenicola"

# row 90
$ws.Cells.Item(90, 1).Value = 75
$ws.Cells.Item(90, 2).Value = "FLAS04"
$ws.Cells.Item(90, 3).Value = "S. elliptica"

# row 91
$ws.Cells.Item(91, 1).Value = 76
$ws.Cells.Item(91, 2).Value = "FLAS05"
$ws.Cells.Item(91, 3).Value = "S. integrifolia"

# row 92
$ws.Cells.Item(92, 1).Value = 77
$ws.Cells.Item(92, 2).Value = "FLAS06"
$ws.Cells.Item(92, 3).Value = "S. multiglandulosa"

# row 93
$ws.Cells.Item(93, 1).Value = 78
$ws.Cells.Item(93, 2).Value = "FLAS07"
$ws.Cells.Item(93, 3).Value = "S. racemosa"

# row 94
$ws.Cells.Item(94, 3).Value = "S.hastifolia"

# row 95
$ws.Cells.Item(95, 3).Value = "S.arenicola"

# row 96
$ws.Cells.Item(96, 3).Value = "S.tournefortii"

# row 97
$ws.Cells.Item(97, 3).Value = "S.baicalensis"

# row 98
$ws.Cells.Item(98, 3).Value = "S.barbata"

# row 99
$ws.Cells.Item(99, 3).Value = "S.indica"

# row 100
$ws.Cells.Item(100, 3).Value = "S.Pekinesis"

# row 101
$ws.Cells.Item(101, 3).Value = "S.dependens"

# row 102
$ws.Cells.Item(102, 3).Value = "S.strigillosa"

# row 103
$ws.Cells.Item(103, 3).Value = "S.insignis"

